$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 117, pushing existing rows 117-153 down to 118-154
$ws.Rows.Item(117).Insert()

# Populate the new row 117 with the new data record
$ws.Cells.Item(117, 1).Value = 1
$ws.Cells.Item(117, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(117, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(117, 4).Value = 44524
$ws.Cells.Item(117, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(117, 5).Value = 15
$ws.Cells.Item(117, 6).Value = "Fruta"
$ws.Cells.Item(117, 7).Value = 100108
$ws.Cells.Item(117, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(117, 9).Value = 100108006
$ws.Cells.Item(117, 10).Value = "Plátano"
$ws.Cells.Item(117, 11).Value = "Barraganete"
$ws.Cells.Item(117, 12).Value = "Primera"
$ws.Cells.Item(117, 13).Value = 120
$ws.Cells.Item(117, 14).Value = 25000
$ws.Cells.Item(117, 15).Value = 26000
$ws.Cells.Item(117, 16).Value = 25500
$ws.Cells.Item(117, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(117, 18).Value = "Ecuador"
$ws.Cells.Item(117, 19).Value = 1275
$ws.Cells.Item(117, 20).Value = 20
